$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3872646.5
$ws.Range("I32").Value = 620
$ws.Range("J32").Value = 5361887.5
$ws.Range("K32").Value = 620
$ws.Range("L32").Value = 5361887.5
$ws.Range("M32").Value = -294
$ws.Range("N32").Value = -5362539.5

$ws.Range("H69").Value = 7855
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 7855
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = ""
$ws.Range("M69").Value = ""
$ws.Range("N69").Value = -25313

$ws.Range("H72").Value = 7855
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 7855
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = ""
$ws.Range("M72").Value = ""
$ws.Range("N72").Value = -79431

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 224.54546
$ws.Range("I5").Value = 140
$ws.Range("J5").Value = 372.5
$ws.Range("K5").Value = 140
$ws.Range("L5").Value = 372.5
$ws.Range("M5").Value = -28
$ws.Range("N5").Value = -596.5

$ws.Range("H63").Value = 17966.666
$ws.Range("I63").Value = 25950
$ws.Range("K63").Value = 25950
$ws.Range("M63").Value = -25264

$ws.Range("H66").Value = 17966.666
$ws.Range("I66").Value = 25950
$ws.Range("K66").Value = 129750
$ws.Range("M66").Value = -126318

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 224.54546
$ws.Range("I4").Value = 140
$ws.Range("J4").Value = 372.5
$ws.Range("K4").Value = 140
$ws.Range("L4").Value = 372.5
$ws.Range("M4").Value = -25
$ws.Range("N4").Value = -602.5

$ws.Range("H82").Value = 11997.637
$ws.Range("I82").Value = 5889.25
$ws.Range("J82").Value = 28286.666
$ws.Range("K82").Value = 5889.25
$ws.Range("L82").Value = 28286.666
$ws.Range("M82").Value = -5506.25
$ws.Range("N82").Value = -29052.666

$ws.Range("H85").Value = 11997.637
$ws.Range("I85").Value = 5889.25
$ws.Range("J85").Value = 28286.666
$ws.Range("K85").Value = 5889.25
$ws.Range("L85").Value = 28286.666
$ws.Range("M85").Value = -4563.25
$ws.Range("N85").Value = -30938.666

$ws.Range("H94").Value = 901.8946999999999
$ws.Range("I94").Value = 766.8570999999999
$ws.Range("J94").Value = 1280
$ws.Range("K94").Value = 766.8570999999999
$ws.Range("L94").Value = 1280
$ws.Range("M94").Value = -315.8570999999999
$ws.Range("N94").Value = -2182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 68000340
$ws.Range("I6").Value = 1500000
$ws.Range("J6").Value = 201001010
$ws.Range("K6").Value = 1500000
$ws.Range("L6").Value = 201001010
$ws.Range("M6").Value = -1499887
$ws.Range("N6").Value = -201001236

$ws.Range("H7").Value = 161.83333
$ws.Range("I7").Value = 189
$ws.Range("J7").Value = 107.5
$ws.Range("K7").Value = 189
$ws.Range("L7").Value = 107.5
$ws.Range("M7").Value = -76
$ws.Range("N7").Value = -333.5

$ws.Range("H17").Value = 3000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = ""

$ws.Range("H22").Value = 384.54544
$ws.Range("I22").Value = 237.5
$ws.Range("J22").Value = 468.57144
$ws.Range("K22").Value = 237.5
$ws.Range("L22").Value = 468.57144
$ws.Range("M22").Value = 112.5
$ws.Range("N22").Value = -1168.57144

$ws.Range("H31").Value = 2504.1462
$ws.Range("I31").Value = 935.63635
$ws.Range("J31").Value = 8974.25
$ws.Range("K31").Value = 935.63635
$ws.Range("L31").Value = 8974.25
$ws.Range("M31").Value = -640.63635
$ws.Range("N31").Value = -9564.25

$ws.Range("H34").Value = 2504.1462
$ws.Range("I34").Value = 935.63635
$ws.Range("J34").Value = 8974.25
$ws.Range("K34").Value = 935.63635
$ws.Range("L34").Value = 8974.25
$ws.Range("M34").Value = -733.63635
$ws.Range("N34").Value = -9378.25

$ws.Range("H38").Value = 7000
$ws.Range("I38").Value = 500
$ws.Range("J38").Value = 13500
$ws.Range("K38").Value = 500
$ws.Range("L38").Value = 13500
$ws.Range("M38").Value = -123
$ws.Range("N38").Value = -14254

$ws.Range("H41").Value = 9958
$ws.Range("I41").Value = 8000
$ws.Range("J41").Value = 12895
$ws.Range("K41").Value = 8000
$ws.Range("L41").Value = 12895
$ws.Range("M41").Value = -7572
$ws.Range("N41").Value = -13751

$ws.Range("H46").Value = 7000
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 13500
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 13500
$ws.Range("M46").Value = -289
$ws.Range("N46").Value = -13922

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = ""
$ws.Range("N51").Value = ""

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = ""
$ws.Range("N59").Value = ""

$ws.Range("H60").Value = 10524.723
$ws.Range("J60").Value = 10949.706
$ws.Range("L60").Value = 10949.706
$ws.Range("N60").Value = -11971.706

$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = ""
$ws.Range("N61").Value = ""

$ws.Range("H68").Value = 32870
$ws.Range("J68").Value = 32870
$ws.Range("L68").Value = 32870
$ws.Range("N68").Value = -34368

$ws.Range("H71").Value = 32870
$ws.Range("J71").Value = 32870
$ws.Range("L71").Value = 98610
$ws.Range("N71").Value = -106098

$ws.Range("H74").Value = 18349
$ws.Range("J74").Value = 18349
$ws.Range("L74").Value = 18349
$ws.Range("N74").Value = -20097

$ws.Range("H77").Value = 18349
$ws.Range("J77").Value = 18349
$ws.Range("L77").Value = 55047
$ws.Range("N77").Value = -63783

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 2000.8334
$ws.Range("I108").Value = 991.6667
$ws.Range("J108").Value = 3010
$ws.Range("K108").Value = 2975.0001
$ws.Range("L108").Value = 9030
$ws.Range("M108").Value = -95.0001000000002
$ws.Range("N108").Value = -14790

$ws.Range("H129").Value = 1300.7949
$ws.Range("I129").Value = 685.1539
$ws.Range("J129").Value = 1608.6154
$ws.Range("K129").Value = 2055.4617
$ws.Range("L129").Value = 4825.8462
$ws.Range("M129").Value = 2944.5383
$ws.Range("N129").Value = -14825.8462

$ws.Range("H131").Value = 1984.0441
$ws.Range("I131").Value = 4646
$ws.Range("J131").Value = 1678.5737
$ws.Range("K131").Value = 13938
$ws.Range("L131").Value = 5035.7211
$ws.Range("M131").Value = -8898
$ws.Range("N131").Value = -15115.7211

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2989.2856
$ws.Range("I126").Value = 2383.3333
$ws.Range("J126").Value = 3443.75
$ws.Range("K126").Value = 7149.999899999999
$ws.Range("L126").Value = 10331.25
$ws.Range("M126").Value = -4679.999899999999
$ws.Range("N126").Value = -15271.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = ""
$ws.Range("N63").Value = ""

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = ""
$ws.Range("N66").Value = ""

$ws.Range("H93").Value = 1633.7858
$ws.Range("I93").Value = 1673.3
$ws.Range("J93").Value = 1535
$ws.Range("K93").Value = 1673.3
$ws.Range("L93").Value = 1535
$ws.Range("M93").Value = -425.3
$ws.Range("N93").Value = -4031

$ws.Range("H98").Value = 26333.334
$ws.Range("J98").Value = 26333.334
$ws.Range("L98").Value = 26333.334
$ws.Range("N98").Value = -32323.334

$ws.Range("H136").Value = 4878
$ws.Range("I136").Value = 3039.6843
$ws.Range("J136").Value = 6932.5884
$ws.Range("K136").Value = 9119.052899999999
$ws.Range("L136").Value = 20797.7652
$ws.Range("M136").Value = -6569.052899999999
$ws.Range("N136").Value = -25897.7652

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 20602580
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 25752976
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 25752976
$ws.Range("M14").Value = -832
$ws.Range("N14").Value = -25753312

$ws.Range("H15").Value = 2220
$ws.Range("J15").Value = 2220
$ws.Range("L15").Value = 2220
$ws.Range("N15").Value = -2796

$ws.Range("H117").Value = 29000
$ws.Range("J117").Value = 29000
$ws.Range("L117").Value = 29000
$ws.Range("N117").Value = -38178
